$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 170.1579
$ws.Range("I8").Value = 80.92857
$ws.Range("J8").Value = 420
$ws.Range("K8").Value = 242.78571
$ws.Range("L8").Value = 1260
$ws.Range("M8").Value = -103.78571
$ws.Range("N8").Value = -1538

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 672.21875
$ws.Range("I15").Value = 672.21875
$ws.Range("K15").Value = 2016.65625
$ws.Range("M15").Value = -1847.65625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1166.194
$ws.Range("J17").Value = 1176.7576
$ws.Range("L17").Value = 3530.2728
$ws.Range("N17").Value = -3866.2728

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 14101.571
$ws.Range("I18").Value = 20427.75
$ws.Range("J18").Value = 5666.6665
$ws.Range("K18").Value = 20427.75
$ws.Range("L18").Value = 5666.6665
$ws.Range("M18").Value = -20143.75
$ws.Range("N18").Value = -6234.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5687.2856
$ws.Range("I43").Value = 5083
$ws.Range("J43").Value = 5929
$ws.Range("K43").Value = 5083
$ws.Range("L43").Value = 5929
$ws.Range("M43").Value = -5014
$ws.Range("N43").Value = -6067

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 100009250
$ws.Range("I69").Value = 5499
$ws.Range("K69").Value = 16497
$ws.Range("M69").Value = -15623

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 100009250
$ws.Range("I72").Value = 5499
$ws.Range("K72").Value = 49491
$ws.Range("M72").Value = -45123

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4408.8096
$ws.Range("I100").Value = 2048.9167
$ws.Range("K100").Value = 2048.9167
$ws.Range("M100").Value = -1507.9167

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 100750
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 100750
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 906750
$ws.Range("N125").Value = -911670
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2074.3
$ws.Range("I132").Value = 1749.36
$ws.Range("K132").Value = 5248.08
$ws.Range("M132").Value = -2718.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 10765.154
$ws.Range("J63").Value = 9995.363
$ws.Range("L63").Value = 9995.363
$ws.Range("N63").Value = -11367.363

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 10765.154
$ws.Range("J66").Value = 9995.363
$ws.Range("L66").Value = 49976.815
$ws.Range("N66").Value = -56840.815

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2571.75
$ws.Range("I122").Value = 2409.8333
$ws.Range("K122").Value = 7229.499899999999
$ws.Range("M122").Value = -4779.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 676
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -326

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 75000
$ws.Range("J111").Value = 75000
$ws.Range("L111").Value = 75000
$ws.Range("N111").Value = -83180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1057.5714
$ws.Range("I22").Value = 1157
$ws.Range("J22").Value = 1007.8571
$ws.Range("K22").Value = 1157
$ws.Range("L22").Value = 1007.8571
$ws.Range("M22").Value = -807
$ws.Range("N22").Value = -1707.8571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3981.2307
$ws.Range("I31").Value = 1322.6428
$ws.Range("K31").Value = 1322.6428
$ws.Range("M31").Value = -1027.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3981.2307
$ws.Range("I34").Value = 1322.6428
$ws.Range("K34").Value = 1322.6428
$ws.Range("M34").Value = -1120.6428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3249
$ws.Range("I99").Value = 3249
$ws.Range("K99").Value = 3249
$ws.Range("M99").Value = -1751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3249
$ws.Range("I126").Value = 3249
$ws.Range("K126").Value = 9747
$ws.Range("M126").Value = -7277

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1408.3429
$ws.Range("I132").Value = 1331.2188
$ws.Range("K132").Value = 3993.6564
$ws.Range("M132").Value = -1463.6564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1990
$ws.Range("I5").Value = 1990
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5970
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -5858
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1990
$ws.Range("I135").Value = 1990
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 17910
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -15375
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7151.857
$ws.Range("I70").Value = 7028.8823
$ws.Range("J70").Value = 7674.5
$ws.Range("K70").Value = 7028.8823
$ws.Range("L70").Value = 7674.5
$ws.Range("M70").Value = -6758.8823
$ws.Range("N70").Value = -8214.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7151.857
$ws.Range("I73").Value = 7028.8823
$ws.Range("J73").Value = 7674.5
$ws.Range("K73").Value = 7028.8823
$ws.Range("L73").Value = 7674.5
$ws.Range("M73").Value = -6092.8823
$ws.Range("N73").Value = -9546.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6726.636
$ws.Range("I102").Value = 4332.3335
$ws.Range("J102").Value = 9599.8
$ws.Range("K102").Value = 4332.3335
$ws.Range("L102").Value = 9599.8
$ws.Range("M102").Value = -2710.3335
$ws.Range("N102").Value = -12843.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2999
$ws.Range("I132").Value = 2999
$ws.Range("K132").Value = 8997
$ws.Range("M132").Value = -6467

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3969.84
$ws.Range("I61").Value = 2735.889
$ws.Range("K61").Value = 2735.889
$ws.Range("M61").Value = -2533.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5241.5386
$ws.Range("I68").Value = 3157
$ws.Range("J68").Value = 6544.375
$ws.Range("K68").Value = 3157
$ws.Range("L68").Value = 6544.375
$ws.Range("M68").Value = -2408
$ws.Range("N68").Value = -8042.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5241.5386
$ws.Range("I71").Value = 3157
$ws.Range("J71").Value = 6544.375
$ws.Range("K71").Value = 15785
$ws.Range("L71").Value = 32721.875
$ws.Range("M71").Value = -12041
$ws.Range("N71").Value = -40209.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3969.84
$ws.Range("I113").Value = 2735.889
$ws.Range("K113").Value = 2735.889
$ws.Range("M113").Value = -565.8890000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3901.44
$ws.Range("I132").Value = 3726.6
$ws.Range("K132").Value = 11179.8
$ws.Range("M132").Value = -8649.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3228.7058
$ws.Range("I122").Value = 1866.3684
$ws.Range("K122").Value = 5599.1052
$ws.Range("M122").Value = -3149.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1806.5
$ws.Range("I132").Value = 1815.1818
$ws.Range("K132").Value = 5445.5454
$ws.Range("M132").Value = -2915.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1526.4546
$ws.Range("I136").Value = 1229.1
$ws.Range("K136").Value = 3687.3
$ws.Range("M136").Value = -1137.3
